$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the national_id value for row 14 (B14): remove the stray trailing
# zero-width space from "079534653344 \u200b" -> "079534653344".
# A leading apostrophe preserves the existing text/quote-prefix style (s=6)
# instead of Excel re-evaluating the value and dropping quotePrefix.
$ws.Cells.Item(14, 2).Value = "'079534653344"

# Update the selected cell/range shown when the workbook is next opened.
$ws.Range("C17").Select()
